$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value to include decimal formatting
$ws.Range("B2").Value = "0.0,0.0"

# Remove rows 3-6 entirely (data for additional fixation locations)
$ws.Range("A3:D6").EntireRow.Delete()
